# Autogenerated on Mon Feb 09 2015 03:30:35 GMT+0000 (Coordinated Universal Time)
#
# Adds a new "Enterprise Size Classification" style table (Number of
# employees / Assets / Turnover, broken out by Micro / Small / Medium /
# Large) to the Summary sheet, between the existing "Value added to the
# economy" block (ending at row 17) and the "Sector Distribution Details"
# block (which was at rows 22/25/26 and needs to shift down to make room).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: push everything from row 19 down by 6 rows. This shifts the
# old "Sector Distribution Details" (22) / "SCRUS" (25) / source note (26)
# rows down to 28 / 31 / 32, matching the target layout, without touching
# their existing content or formatting.
$ws.Rows("19:24").Insert()

# New table header row (bold, like the other "Micro / SMEs / MSMEs"-style
# header rows already on the sheet).
$ws.Cells.Item(19, 2).Value = "Number of employees"
$ws.Cells.Item(19, 2).Font.Bold = $true
$ws.Cells.Item(19, 3).Value = "Assets (local currency, unless noted otherwise)"
$ws.Cells.Item(19, 3).Font.Bold = $true
$ws.Cells.Item(19, 4).Value = "Turnover (local currency, unless noted otherwise)"
$ws.Cells.Item(19, 4).Font.Bold = $true

# Row labels for the new table.
$ws.Cells.Item(20, 1).Value = "Micro"
$ws.Cells.Item(21, 1).Value = "Small"
$ws.Cells.Item(22, 1).Value = "Medium"
$ws.Cells.Item(23, 1).Value = "Large"

# The data columns (employees / assets / turnover) have no values yet for
# this country, so each cell is present but holds an empty string.
for ($r = 20; $r -le 23; $r++) {
  for ($c = 2; $c -le 4; $c++) {
    $ws.Cells.Item($r, $c).Formula = "=""""" 
  }
}
